$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new "Status" column ---
$ws.Range("D1").Value = "สถานะ"

# --- Row 2 (item 1 - Select Class): ยังไม่ได้แก้ไข (not fixed yet) -> red ---
$ws.Range("D2").Value = "ยังไม่ได้แก้ไข"
$ws.Range("D2").Interior.Color = 8420607

# --- Row 3 (item 2 - Select Lane): ยังไม่ได้แก้ไข -> red ---
$ws.Range("D3").Value = "ยังไม่ได้แก้ไข"
$ws.Range("D3").Interior.Color = 8420607

# --- Row 4 (item 3 - Select Support Item): แก้ไขแล้ว (fixed) -> green/theme ---
$ws.Range("D4").Value = "แก้ไขแล้ว"
$ws.Range("D4").Interior.ThemeColor = 10

# --- Row 5 (item 4 - Select Farm Item): แก้ไขแล้ว -> green/theme ---
$ws.Range("D5").Value = "แก้ไขแล้ว"
$ws.Range("D5").Interior.ThemeColor = 10

# --- Row 6 (new item 5 - Select Support Item, Support Lane note): red ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Select Support Item"
$ws.Range("C7").Value = "แก้ไขให้แสดงเมื่อเลือก Farm Lane ใน Select Lane"
$ws.Range("C6").Value = "แก้ไขให้แสดงเมื่อเลือก Support Lane ใน Select Lane"
$ws.Range("D6").Value = "ยังไม่ได้แก้ไข"
$ws.Range("D6").Interior.Color = 8420607

# --- Row 7 (new item 6 - Select Farm Item, Farm Lane note): red ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Select Farm Item"
$ws.Range("D7").Value = "ยังไม่ได้แก้ไข"
$ws.Range("D7").Interior.Color = 8420607

# --- Rows 8-18: trailing numbered rows, column A only ---
for ($n = 7; $n -le 17; $n++) {
    $r = $n + 1
    $ws.Cells.Item($r, 1).Value = $n
}

# --- Row heights / column widths ---
$ws.Rows.Item(2).RowHeight = 33.65
$ws.Columns.Item(3).ColumnWidth = 65.666666666666667
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668

# --- Selection ---
$ws.Range("F5").Select()
